$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "panduan lelang"

# Reset column B's width to the sheet's standard/default width before removing
# column C, so that the OOXML <cols> list stays well-formed (avoid a buggy
# min/max shift that occurs when deleting a column that still carries an
# explicit custom width).
$ws.Columns("B").ColumnWidth = $ws.Columns("D").ColumnWidth

# Replace the "judul" column header and its values with the "expected" column
# (which used to live in column C) and drop the redundant per-row result
# column entirely.
$ws.Range("B1").Value = "expected"
$ws.Range("B2").Value = "pass"
$ws.Range("B3").Value = "pass"
$ws.Range("B4").Value = "pass"
$ws.Range("B5").Value = "pass"
$ws.Range("B6").Value = "pass"

# Drop the now-duplicated column C (old "expected" column).
$ws.Columns("C").Delete()

# Drop the last "semua" summary row.
$ws.Rows("6").Delete()
